$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a weekly-refreshed price log ("Fruta / hortaliza, semanal").
# Two new weekly rows of data were inserted into the existing table, pushing
# the rows that used to follow them down by one each. We reproduce this with
# two row inserts (each insert shifts everything at/after that row down by
# one), then populate the two freshly-inserted blank rows with their data.
# ---------------------------------------------------------------------------

# 1) Insert the first new row at (current) row 352 - shifts old rows 352..467
#    down to 353..468.
$ws.Rows.Item(352).Insert()

# 2) Insert the second new row at (current) row 464 - this position now holds
#    what used to be old row 463 (after the first shift); inserting here
#    shifts it (and everything after) down one more, to 465..469.
$ws.Rows.Item(464).Insert()

# ---------------------------------------------------------------------------
# Fill in the new row 352 (copy of the static columns from its neighbour,
# plus the new measurement values).
# ---------------------------------------------------------------------------
$ws.Cells.Item(352, 1).Value = 4
$ws.Cells.Item(352, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(352, 3).Value = "Los Lagos"
$ws.Cells.Item(352, 4).Value = 45120
$ws.Cells.Item(352, 5).Value = 10
$ws.Cells.Item(352, 6).Value = 100112043
$ws.Cells.Item(352, 7).Value = "Pepino ensalada"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 150
$ws.Cells.Item(352, 11).Value = 17000
$ws.Cells.Item(352, 12).Value = 18000
$ws.Cells.Item(352, 13).Value = 17533
$ws.Cells.Item(352, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(352, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(352, 16).Value = 292
$ws.Cells.Item(352, 17).Value = 60
$ws.Cells.Item(352, 18).Value = "Hortaliza"

# ---------------------------------------------------------------------------
# Fill in the new row 464 (copy of the static columns from its neighbour,
# plus the new measurement values).
# ---------------------------------------------------------------------------
$ws.Cells.Item(464, 1).Value = 4
$ws.Cells.Item(464, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(464, 3).Value = "Los Lagos"
$ws.Cells.Item(464, 4).Value = 45121
$ws.Cells.Item(464, 5).Value = 10
$ws.Cells.Item(464, 6).Value = 100112043
$ws.Cells.Item(464, 7).Value = "Pepino ensalada"
$ws.Cells.Item(464, 8).Value = "Sin especificar"
$ws.Cells.Item(464, 9).Value = "Primera"
$ws.Cells.Item(464, 10).Value = 400
$ws.Cells.Item(464, 11).Value = 16500
$ws.Cells.Item(464, 12).Value = 17000
$ws.Cells.Item(464, 13).Value = 16750
$ws.Cells.Item(464, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(464, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(464, 16).Value = 279
$ws.Cells.Item(464, 17).Value = 60
$ws.Cells.Item(464, 18).Value = "Hortaliza"

# Keep the date columns formatted the same way as the rest of column D.
$ws.Cells.Item(352, 4).Style = $ws.Cells.Item(353, 4).Style
$ws.Cells.Item(464, 4).Style = $ws.Cells.Item(465, 4).Style
